# Applies the "Improved program. Runs autonomously well but not perfect" edit
# to the Snake Game "Directions" decision table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data changes -----------------------------------------------------
# Row 3 ("Left" rule): the "Front" sensor flag changes from 1 to 0.
$ws.Range("C3").Value = 0

# Row 5: the resulting action changes from "Front" to "Right", and a new
# comment cell "#Front" is added next to it (column H).
$ws.Range("G5").Value = "Right"
$ws.Range("H5").Value = "#Front"

# Row 15: the resulting action changes from "Front" to "Left", and a new
# comment cell "#Front" is added next to it (column H).
$ws.Range("G15").Value = "Left"
$ws.Range("H15").Value = "#Front"

# --- Selection / active cell ------------------------------------------
# The saved selection moves from H16 to J13.
$ws.Range("J13").Select()
